# Add the new recipient row (glevon111@gmail.com) to the sheet, with its
# mailto hyperlink, mirroring the formatting used by the other hyperlink
# rows already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7 holds the new e-mail address.
$newCell = $ws.Range("A7")
$newCell.Value = "glevon111@gmail.com"

# Turn it into a live mailto: hyperlink, like the other addresses.
$ws.Hyperlinks.Add($newCell, "mailto:glevon111@gmail.com")

# Hyperlinks.Add() stamps its own cell style; re-apply the same
# formatting already used by the existing hyperlink cells (e.g. A5) so
# the new cell matches them exactly.
$ws.Range("A5").Copy()
$newCell.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection down to A8, just like Excel does after data
# entry in A7.
$ws.Range("A8").Select() | Out-Null
